$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 179: "Al Hmanah" / "الحمنة" -> "Urayarah" / "عريعرة", update coordinates and area/region
$ws.Range("A179").Value = "Urayarah"
$ws.Range("B179").Value = "Urayarah"
$ws.Range("C179").Value = "عريعرة"
$ws.Range("D179").Value = 25.980965999999999
$ws.Range("E179").Value = 48.849434000000002
$ws.Range("F179").Value = "المنطقة الشرقية"
$ws.Range("G179").Value = "شرق المملكة"

# Row 180: "Al Qahma" / "القحمة" -> "Al Mossam" / "الموسم", update coordinates and area
$ws.Range("A180").Value = "Al Mossam"
$ws.Range("B180").Value = "Al Mossam"
$ws.Range("C180").Value = "الموسم"
$ws.Range("D180").Value = 16.418261999999999
$ws.Range("E180").Value = 42.824491999999999
$ws.Range("F180").Value = "منطقة جازان"
